$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 0.03742607606442153
$ws.Range("C2").Value = 0.01735433991663465
$ws.Range("D2").Value = 0.03506375804688142
$ws.Range("E2").Value = 0.008118037228071846
$ws.Range("F2").Value = 0.03682596592847593
$ws.Range("G2").Value = 0.0200845498731792
$ws.Range("B3").Value = 0.02806055989938596
$ws.Range("C3").Value = 0.0287711510192422
$ws.Range("D3").Value = 0.04918539443472066
$ws.Range("E3").Value = 0.04911947615978158
$ws.Range("F3").Value = 0.01928526081147312
$ws.Range("G3").Value = 0.01986233663867399
$ws.Range("B4").Value = 0.03646409467464339
$ws.Range("C4").Value = 0.03580779814367356
$ws.Range("D4").Value = 0.02011035498457307
$ws.Range("E4").Value = 0.02754968594147661
$ws.Range("F4").Value = 0.01064380226166979
$ws.Range("G4").Value = 0.01048723923400685
$ws.Range("B5").Value = 0.02118225111743641
$ws.Range("C5").Value = 0.02103491568459112
$ws.Range("D5").Value = 0.01033833192424172
$ws.Range("E5").Value = 0.01145276069530116
$ws.Range("F5").Value = 0.004932906579649861
$ws.Range("G5").Value = 0.00466256026706771
$ws.Range("B6").Value = 0.1452262427004578
$ws.Range("C6").Value = 0.03881545611889942
$ws.Range("D6").Value = 0.0815887606899436
$ws.Range("E6").Value = 0.02910867774809549
$ws.Range("F6").Value = 0.1165753023852861
$ws.Range("G6").Value = 0.04179805946082008
$ws.Range("B7").Value = 0.0399857475402855
$ws.Range("C7").Value = 0.02953529607767487
$ws.Range("D7").Value = 0.01651853452041685
$ws.Range("E7").Value = 0.01555199674721331
$ws.Range("F7").Value = 0.02569016821914447
$ws.Range("G7").Value = 0.01949593085650963
$ws.Range("B8").Value = 0.8637590537089674
$ws.Range("C8").Value = 0.8861858005738782
$ws.Range("D8").Value = 0.9314824565938706
$ws.Range("E8").Value = 0.9385707458446433
$ws.Range("F8").Value = 0.8905505363486729
$ws.Range("G8").Value = 0.9051044498576281
$ws.Range("B9").Value = 0.04121757406691675
$ws.Range("C9").Value = 0.04200305340334829
$ws.Range("D9").Value = 0.01340434543259565
$ws.Range("E9").Value = 0.01338077594820928
$ws.Range("F9").Value = 0.02398291086272598
$ws.Range("G9").Value = 0.02193740716661464
$ws.Range("C10").Value = 0.1602595465005063
$ws.Range("E10").Value = 0.08749139758189517
$ws.Range("G10").Value = 0.1196837274874851
$ws.Range("C11").Value = 0.04229280253821115
$ws.Range("E11").Value = 0.02344653717128905
$ws.Range("G11").Value = 0.02840974511261153
$ws.Range("B12").Value = 4.291113069285608
$ws.Range("C12").Value = 4.329020436629127
$ws.Range("D12").Value = 3.872445733506573
$ws.Range("E12").Value = 3.891870898365699
$ws.Range("F12").Value = 6.698973456559248
$ws.Range("G12").Value = 7.095796123900807
$ws.Range("B13").Value = 0.4186904290563454
$ws.Range("C13").Value = 0.4159186518916121
$ws.Range("D13").Value = 0.332751578697806
$ws.Range("E13").Value = 0.3363378721560813
$ws.Range("F13").Value = 0.7775384444035005
$ws.Range("G13").Value = 0.8625573269246137
$ws.Range("B14").Value = -4833.862741016288
$ws.Range("C14").Value = -4818.973508420826
$ws.Range("D14").Value = -6441.222436456836
$ws.Range("E14").Value = -6433.000686815352
$ws.Range("F14").Value = -3961.195665966567
$ws.Range("G14").Value = -3949.595743188088
$ws.Range("B15").Value = 9675.725482032576
$ws.Range("C15").Value = 9647.947016841652
$ws.Range("D15").Value = 12890.44487291367
$ws.Range("E15").Value = 12876.0013736307
$ws.Range("F15").Value = 7930.391331933134
$ws.Range("G15").Value = 7909.191486376176
$ws.Range("B16").Value = 9699.021666076002
$ws.Range("C16").Value = 9677.067246895933
$ws.Range("D16").Value = 12913.7410569571
$ws.Range("E16").Value = 12905.12160368499
$ws.Range("F16").Value = 7953.68751597656
$ws.Range("G16").Value = 7938.311716430458

Write-Output "Done updating cells"
